$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the placeholder picture/drawing (rId1) that was anchored at J4.
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# Row 4 height shrinks now that the picture is gone.
$ws.Rows("4").RowHeight = 60

# J4 now carries descriptive text instead of the (removed) picture.
$ws.Range("J4").Value = "直鋼筋 安#3-390x40" + [char]10 + "長度: 390cm"

# Refresh the generated-at timestamp footer.
$ws.Range("A9").Value = "生成時間：2025-10-06 08:08:06 | 圖示功能：暫時停用，等 assets/materials 圖片準備好時再實作"
